$d = $word.ActiveDocument

# The paragraph's final sentence ("We will derive ... This chapter shows")
# is being reworded. The "_GoBack" bookmark currently sits between the
# "...European scale. " run and the "We will derive..." run; in the edited
# text it needs to end up right after "(see Appendix 2b and 2c)" and right
# before the (now merged) ". This chapter shows..." sentence. To avoid the
# bookmark being silently dropped by a replace that spans across it, first
# rewrite only the text that comes after the bookmark, then relocate the
# bookmark to its new position with Bookmarks.Add (which moves a bookmark
# that already exists under that name).

$old = "We will derive the green space data from UA, as well as the information on the cities inhabitants. From OSM we will derive the network and the buildings. We will estimate building and park entry points based on both datasets. This chapter shows"
$new = "We derive the green space data as well as the information on the city" + [char]0x2019 + "s inhabitants from UA. From OSM we derive the network and the buildings. We estimate the locations of building and park entry points based on both datasets (see Appendix 2b and 2c). This chapter shows"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $false, 1, $false, $new, 2)

# Relocate the "_GoBack" bookmark to sit right after "(see Appendix 2b and 2c)"
$r = $d.Content
$r.Find.Execute("(see Appendix 2b and 2c)", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
